# Add the new metadata record (row 2) to Sheet1, matching columns:
#   A=identifier  B=alternativeIdentifiers  C=title  D=date_s
#   E=levelOfDescription  F=extentAndMedium  G=notes  H=file_path

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style the new row's cells (A, C, D, E, F, G, H -- B2 is left untouched)
# with the same body font used elsewhere in the sheet (Calibri 10, theme text color).
$bodyCells = @("A2", "C2", "D2", "E2", "F2", "G2", "H2")
foreach ($addr in $bodyCells) {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 10
    $cell.Font.ThemeColor = 1
}

$ws.Range("A2").Value = "MCH185"
$ws.Range("C2").Value = "OPEN APARTHEID PRISONS, ANC BEYOND APARTHEID"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

$ws.Range("A2:I2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
